$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "'33"
$ws.Cells.Item(2,2).Value = "'27"
$ws.Cells.Item(2,4).Value = 1006
$ws.Cells.Item(2,5).Value = 66
$ws.Cells.Item(2,6).Value = 71

# Row 3
$ws.Cells.Item(3,1).Value = "'33"
$ws.Cells.Item(3,4).Value = 1007
$ws.Cells.Item(3,5).Value = 25
$ws.Cells.Item(3,6).Value = 69

# Row 4
$ws.Cells.Item(4,2).Value = "'29"
$ws.Cells.Item(4,4).Value = 1008
$ws.Cells.Item(4,5).Value = 40
$ws.Cells.Item(4,6).Value = 71

# Row 5
$ws.Cells.Item(5,1).Value = "'33"
$ws.Cells.Item(5,4).Value = 1007
$ws.Cells.Item(5,6).Value = 74

# Row 6
$ws.Cells.Item(6,2).Value = "'28"
$ws.Cells.Item(6,4).Value = 1006
$ws.Cells.Item(6,5).Value = 71
$ws.Cells.Item(6,6).Value = 79

# Row 7
$ws.Cells.Item(7,1).Value = "'32"
$ws.Cells.Item(7,4).Value = 1008
$ws.Cells.Item(7,5).Value = 61
$ws.Cells.Item(7,6).Value = 79

# Row 8
$ws.Cells.Item(8,1).Value = "'34"
$ws.Cells.Item(8,2).Value = "'25"
$ws.Cells.Item(8,4).Value = 1009
$ws.Cells.Item(8,5).Value = 38
$ws.Cells.Item(8,6).Value = 73
